$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04760794252374
$ws.Range("D2").Value = 1.056221041139173
$ws.Range("E2").Value = 1.061090968700482
$ws.Range("F2").Value = 1.068201649286932
$ws.Range("I2").Value = 1.047975993275452
$ws.Range("J2").Value = 1.052656106692502
$ws.Range("K2").Value = 1.058959347971665
$ws.Range("L2").Value = 1.063815972625804
$ws.Range("M2").Value = 1.070907465045288
$ws.Range("N2").Value = 1.021392693958611

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.04848370412149
$ws.Range("D3").Value = 1.056935287091577
$ws.Range("E3").Value = 1.06193200531932
$ws.Range("F3").Value = 1.069072287641119
$ws.Range("I3").Value = 1.04822298765742
$ws.Range("J3").Value = 1.053180772434348
$ws.Range("K3").Value = 1.05948744313225
$ws.Range("L3").Value = 1.064471495570295
$ws.Range("M3").Value = 1.071593899813898
$ws.Range("N3").Value = 1.021568660645228

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.049050961230855
$ws.Range("D4").Value = 1.057397909138205
$ws.Range("E4").Value = 1.06247712349885
$ws.Range("F4").Value = 1.069636586161568
$ws.Range("I4").Value = 1.048381754242625
$ws.Range("J4").Value = 1.053520177255426
$ws.Range("K4").Value = 1.059828938472189
$ws.Range("L4").Value = 1.06489592408991
$ws.Range("M4").Value = 1.072038361319626
$ws.Range("N4").Value = 1.021682449418525

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.049289573706918
$ws.Range("D5").Value = 1.057592503230044
$ws.Range("E5").Value = 1.06270650760954
$ws.Range("F5").Value = 1.069874039872678
$ws.Range("I5").Value = 1.048448246317252
$ws.Range("J5").Value = 1.05366284055673
$ws.Range("K5").Value = 1.059972449771577
$ws.Range("L5").Value = 1.065074414984182
$ws.Range("M5").Value = 1.072225281422237
$ws.Range("N5").Value = 1.02173026819002

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.049329645817521
$ws.Range("D6").Value = 1.057625182706503
$ws.Range("E6").Value = 1.062745034852313
$ws.Range("F6").Value = 1.069913922398155
$ws.Range("I6").Value = 1.0484593957478
$ws.Range("J6").Value = 1.053686792999939
$ws.Range("K6").Value = 1.059996542785587
$ws.Range("L6").Value = 1.065104387934982
$ws.Range("M6").Value = 1.072256670102563
$ws.Range("N6").Value = 1.021738296101947

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04905414904363
$ws.Range("D7").Value = 1.057400508892553
$ws.Range("E7").Value = 1.062480187692314
$ws.Range("F7").Value = 1.069639758157186
$ws.Range("I7").Value = 1.048382643709604
$ws.Range("J7").Value = 1.05352208361927
$ws.Range("K7").Value = 1.05983085628837
$ws.Range("L7").Value = 1.064898308855499
$ws.Range("M7").Value = 1.072040858687956
$ws.Range("N7").Value = 1.021683088446691

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047903789121998
$ws.Range("D8").Value = 1.056462328247826
$ws.Range("E8").Value = 1.061375011245237
$ws.Range("F8").Value = 1.06849569045597
$ws.Range("I8").Value = 1.04805968423316
$ws.Range("J8").Value = 1.05283343765465
$ws.Range("K8").Value = 1.059137864488508
$ws.Range("L8").Value = 1.064037454707204
$ws.Range("M8").Value = 1.071139387466158
$ws.Range("N8").Value = 1.021452177661503

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.045881220601311
$ws.Range("D9").Value = 1.05481270871235
$ws.Range("E9").Value = 1.059434604553497
$ws.Range("F9").Value = 1.066486953985246
$ws.Range("I9").Value = 1.047482539601514
$ws.Range("J9").Value = 1.051619326002321
$ws.Range("K9").Value = 1.05791511405284
$ws.Range("L9").Value = 1.062522584085437
$ws.Range("M9").Value = 1.06955318478499
$ws.Range("N9").Value = 1.021044740657039

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.044535961998785
$ws.Range("D10").Value = 1.053715464489851
$ws.Range("E10").Value = 1.058145846466083
$ws.Range("F10").Value = 1.065152781134358
$ws.Range("I10").Value = 1.047092410014751
$ws.Range("J10").Value = 1.050809565070384
$ws.Range("K10").Value = 1.057098939362692
$ws.Range("L10").Value = 1.061514142772185
$ws.Range("M10").Value = 1.06849735475917
$ws.Range("N10").Value = 1.020772777687282

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.043954208336187
$ws.Range("D11").Value = 1.053240959882479
$ws.Range("E11").Value = 1.05758897108425
$ws.Range("F11").Value = 1.064576274096121
$ws.Range("I11").Value = 1.046922216159624
$ws.Range("J11").Value = 1.0504588605781
$ws.Range("K11").Value = 1.056745303426813
$ws.Range("L11").Value = 1.061077843532347
$ws.Range("M11").Value = 1.068040575973497
$ws.Range("N11").Value = 1.020654940062849

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.043738233563501
$ws.Range("D12").Value = 1.05306480096789
$ws.Range("E12").Value = 1.057382299596038
$ws.Range("F12").Value = 1.064362315578209
$ws.Range("I12").Value = 1.046858809242186
$ws.Range("J12").Value = 1.050328583567737
$ws.Range("K12").Value = 1.05661391445737
$ws.Range("L12").Value = 1.060915838450102
$ws.Range("M12").Value = 1.067870969971796
$ws.Range("N12").Value = 1.020611158970483

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043784555683692
$ws.Range("D13").Value = 1.053102583407656
$ws.Range("E13").Value = 1.057426623302264
$ws.Range("F13").Value = 1.064408202148122
$ws.Range("I13").Value = 1.046872418797497
$ws.Range("J13").Value = 1.050356528850961
$ws.Range("K13").Value = 1.05664209930315
$ws.Range("L13").Value = 1.060950586538622
$ws.Range("M13").Value = 1.067907348206305
$ws.Range("N13").Value = 1.020620550651856

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.043936353451261
$ws.Range("D14").Value = 1.053226396620246
$ws.Range("E14").Value = 1.05757188393736
$ws.Range("F14").Value = 1.064558584497956
$ws.Range("I14").Value = 1.046916978782301
$ws.Range("J14").Value = 1.050448092031844
$ws.Range("K14").Value = 1.056734443436232
$ws.Range("L14").Value = 1.061064450993531
$ws.Range("M14").Value = 1.068026555015527
$ws.Range("N14").Value = 1.020651321324935

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.044029896273165
$ws.Range("D15").Value = 1.05330269442992
$ws.Range("E15").Value = 1.0576614072917
$ws.Range("F15").Value = 1.06465126417929
$ws.Range("I15").Value = 1.046944408592027
$ws.Range("J15").Value = 1.05050450587303
$ws.Range("K15").Value = 1.056791335390616
$ws.Range("L15").Value = 1.061134614089539
$ws.Range("M15").Value = 1.068100010522309
$ws.Range("N15").Value = 1.020670278714275

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.044574587043344
$ws.Range("D16").Value = 1.053746968779545
$ws.Range("E16").Value = 1.0581828291527
$ws.Range("F16").Value = 1.065191067413713
$ws.Range("I16").Value = 1.047103678627472
$ws.Range("J16").Value = 1.050832838741423
$ws.Range("K16").Value = 1.057122404341336
$ws.Range("L16").Value = 1.06154310627451
$ws.Range("M16").Value = 1.068527678272269
$ws.Range("N16").Value = 1.020780596622014

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.044916459565278
$ws.Range("D17").Value = 1.054025814670334
$ws.Range("E17").Value = 1.058510216541843
$ws.Range("F17").Value = 1.065529993884987
$ws.Range("I17").Value = 1.047203246309681
$ws.Range("J17").Value = 1.05103877464668
$ws.Range("K17").Value = 1.057330015281929
$ws.Range("L17").Value = 1.061799440684347
$ws.Range("M17").Value = 1.068796051790232
$ws.Range("N17").Value = 1.020849776134779

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.045115940502437
$ws.Range("D18").Value = 1.054188519445208
$ws.Range("E18").Value = 1.058701288392082
$ws.Range("F18").Value = 1.06572779949014
$ws.Range("I18").Value = 1.047261200334188
$ws.Range("J18").Value = 1.051158886486391
$ws.Range("K18").Value = 1.057451089253552
$ws.Range("L18").Value = 1.061948991089864
$ws.Range("M18").Value = 1.068952628265939
$ws.Range("N18").Value = 1.020890120024501

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.045183970597164
$ws.Range("D19").Value = 1.05424400747182
$ws.Range("E19").Value = 1.058766457954653
$ws.Range("F19").Value = 1.06579526563103
$ws.Range("I19").Value = 1.047280940422863
$ws.Range("J19").Value = 1.051199840281644
$ws.Range("K19").Value = 1.057492368583155
$ws.Range("L19").Value = 1.061999989786678
$ws.Range("M19").Value = 1.069006023323938
$ws.Range("N19").Value = 1.020903874988126

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.044879772371918
$ws.Range("D20").Value = 1.053995891063146
$ws.Range("E20").Value = 1.058475079342481
$ws.Range("F20").Value = 1.065493618319612
$ws.Range("I20").Value = 1.047192576265809
$ws.Range("J20").Value = 1.051016680397868
$ws.Range("K20").Value = 1.057307742857001
$ws.Range("L20").Value = 1.061771934809929
$ws.Range("M20").Value = 1.068767253833437
$ws.Range("N20").Value = 1.020842354584756

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.043891649621541
$ws.Range("D21").Value = 1.053189934128249
$ws.Range("E21").Value = 1.05752910340569
$ws.Range("F21").Value = 1.06451429560939
$ws.Range("I21").Value = 1.046903862198246
$ws.Range("J21").Value = 1.050421129216883
$ws.Range("K21").Value = 1.056707251287056
$ws.Range("L21").Value = 1.061030919195546
$ws.Range("M21").Value = 1.067991449858484
$ws.Range("N21").Value = 1.0206422604274

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.043271040457426
$ws.Range("D22").Value = 1.052683737206308
$ws.Range("E22").Value = 1.056935354695688
$ws.Range("F22").Value = 1.063899609652498
$ws.Range("I22").Value = 1.046721241097659
$ws.Range("J22").Value = 1.050046626740025
$ws.Range("K22").Value = 1.056329509038512
$ws.Range("L22").Value = 1.060565337285217
$ws.Range("M22").Value = 1.067504030429796
$ws.Range("N22").Value = 1.020516390145042

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.043599973809956
$ws.Range("D23").Value = 1.052952030010497
$ws.Range("E23").Value = 1.05725001448562
$ws.Range("F23").Value = 1.064225365889819
$ws.Range("I23").Value = 1.046818155555677
$ws.Range("J23").Value = 1.050245162488127
$ws.Range("K23").Value = 1.056529774882516
$ws.Range("L23").Value = 1.060812119877681
$ws.Range("M23").Value = 1.067762386132666
$ws.Range("N23").Value = 1.020583122199759

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044896349511622
$ws.Range("D24").Value = 1.054009412072065
$ws.Range("E24").Value = 1.058490955983554
$ws.Range("F24").Value = 1.065510054515088
$ws.Range("I24").Value = 1.047197397976757
$ws.Range("J24").Value = 1.051026663859911
$ws.Range("K24").Value = 1.057317806875172
$ws.Range("L24").Value = 1.061784363422991
$ws.Range("M24").Value = 1.068780266270984
$ws.Range("N24").Value = 1.020845708086516

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.046403559387912
$ws.Range("D25").Value = 1.055238741112261
$ws.Range("E25").Value = 1.059935399508989
$ws.Range("F25").Value = 1.067005389640856
$ws.Range("I25").Value = 1.0476326945101
$ws.Range("J25").Value = 1.051933270346727
$ws.Range("K25").Value = 1.058231407460772
$ws.Range("L25").Value = 1.062913961112184
$ws.Range("M25").Value = 1.069962973919996
$ws.Range("N25").Value = 1.021150134462304

